# This script applies a series of text-normalization edits to the
# document: Word's spell/grammar-checker had split several phrases into
# multiple runs wrapped in <w:proofErr .../> markers (spellStart/spellEnd,
# gramStart/gramEnd). The author accepted/cleared all those proofing
# marks, which causes Word to re-merge the adjacent runs back into a
# single run per paragraph. We reproduce that by re-typing each affected
# phrase via Find & Replace (same text in, same text out) - this is
# enough to make Word normalize the runs and drop the proofErr markers.
# In addition, one paragraph gets a brand-new " 1/3 hora" run appended
# after "Esfuerzo real:" (the second "Testing" section).

$d = $word.ActiveDocument

function Retype-Text($text) {
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $text, 2) | Out-Null
}

# --- Names on the cover page -------------------------------------------------
Retype-Text "Agustina Disiot 221025"

# --- Recurring "Requirements Definition" / "Balsamiq Prototyping" / --------
# --- "Testing" section headers (each occurs 3 times in the document) ------
Retype-Text "Requirements Definition"
Retype-Text "Balsamiq Prototyping"
Retype-Text "Testing"

# --- Acceptance-criteria narrative text -------------------------------------
Retype-Text "Para poder hacer una predicción de cuantos casos nuevos pueden haber y "
Retype-Text "Cuando se clickea en la ventana de estadisticas"
Retype-Text "Entonces el sistema despliega la ventana estadisticas que incluye una grafica con los casos actuales y totales de Covid-19"
Retype-Text "Agregar nuevas user stories en la iteración 3"
Retype-Text "Como ciudadano Uruguayo"
Retype-Text "Para acceder a las funcionalidades que tengo permitidas una vez logueado."
Retype-Text "Quiero poder se informado si estuve expuesto al coronavirus."

# --- New "esfuerzo real" value for the Testing task of the second HU -------
# The phrase "Esfuerzo estimado: 1/3 horas ideales" occurs twice in the
# document; the paragraph that needs the new text follows the *second*
# occurrence, so skip past the first hit before searching.
$range1 = $d.Content
$range1.Find.ClearFormatting()
$range1.Find.Forward = $true
$range1.Find.Wrap = 0
$range1.Find.Execute("Esfuerzo estimado: 1/3 horas ideales") | Out-Null

$range2 = $d.Range($range1.End, $d.Content.End)
$range2.Find.ClearFormatting()
$range2.Find.Forward = $true
$range2.Find.Wrap = 0
$range2.Find.Execute("Esfuerzo estimado: 1/3 horas ideales") | Out-Null

if ($range2.Find.Found) {
    $para = $range2.Paragraphs(1).Next()
    $target = $para.Range
    $insertPoint = $d.Range($target.End - 1, $target.End - 1)
    $insertPoint.InsertAfter(" 1/3 hora")
}
